$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.705819096474256
$ws.Range("C2").Value = 5.886732840644895
$ws.Range("D2").Value = 4.698982350471384
$ws.Range("E2").Value = 16.47109017746577
$ws.Range("F2").Value = 23.8043630212011
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("K2").Value = 8.882731979470774
$ws.Range("N2").Value = 17.87599817921527
$ws.Range("O2").Value = 21.24733676098599
$ws.Range("B3").Value = 9.374954354194008
$ws.Range("C3").Value = 5.706250316016121
$ws.Range("D3").Value = 4.652626964583883
$ws.Range("E3").Value = 15.53926467437217
$ws.Range("F3").Value = 23.79061507511111
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("K3").Value = 8.647363175537642
$ws.Range("N3").Value = 17.93844410697893
$ws.Range("O3").Value = 21.2924919510841
$ws.Range("B4").Value = 9.167637212089963
$ws.Range("C4").Value = 5.591466566721013
$ws.Range("D4").Value = 4.623474135509062
$ws.Range("E4").Value = 14.9424130199795
$ws.Range("F4").Value = 23.78909603260898
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("K4").Value = 8.501173258586117
$ws.Range("N4").Value = 17.97850918578589
$ws.Range("O4").Value = 21.32540428907016
$ws.Range("B5").Value = 9.082245421718524
$ws.Range("C5").Value = 5.543735994404134
$ws.Range("D5").Value = 4.611425863154473
$ws.Range("E5").Value = 14.69324880114321
$ws.Range("F5").Value = 23.79021741609503
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("K5").Value = 8.441271340075417
$ws.Range("N5").Value = 17.99527072555968
$ws.Range("O5").Value = 21.34011648575572
$ws.Range("B6").Value = 9.068015606255418
$ws.Range("C6").Value = 5.535753985052296
$ws.Range("D6").Value = 4.609415252793566
$ws.Range("E6").Value = 14.65152531950265
$ws.Range("F6").Value = 23.79050870178548
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("K6").Value = 8.431307628284383
$ws.Range("N6").Value = 17.99808026285974
$ws.Range("O6").Value = 21.34263783351084
$ws.Range("B7").Value = 9.166489074324344
$ws.Range("C7").Value = 5.590826664655144
$ws.Range("D7").Value = 4.62331232243314
$ws.Range("E7").Value = 14.9390763681942
$ws.Range("F7").Value = 23.78910411059337
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("K7").Value = 8.500366604492665
$ws.Range("N7").Value = 17.97873347583927
$ws.Range("O7").Value = 21.32559744439217
$ws.Range("B8").Value = 9.592681041092508
$ws.Range("C8").Value = 5.825349865947812
$ws.Range("D8").Value = 4.683145492745095
$ws.Range("E8").Value = 16.15507130896395
$ws.Range("F8").Value = 23.79818628822194
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("K8").Value = 8.801977782593211
$ws.Range("N8").Value = 17.89717278371593
$ws.Range("O8").Value = 21.26182735092079
$ws.Range("B9").Value = 10.38989363547242
$ws.Range("C9").Value = 6.252119765099915
$ws.Range("D9").Value = 4.79477129594236
$ws.Range("E9").Value = 18.41137254896235
$ws.Range("F9").Value = 23.87087055002837
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("K9").Value = 9.376368424807087
$ws.Range("N9").Value = 17.75083767219069
$ws.Range("O9").Value = 21.17810115012547
$ws.Range("B10").Value = 10.94556887747802
$ws.Range("C10").Value = 6.5435206323015
$ws.Range("D10").Value = 4.873007789175806
$ws.Range("E10").Value = 20.04649757973483
$ws.Range("F10").Value = 23.95756516882097
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("K10").Value = 9.78327021208807
$ws.Range("N10").Value = 17.65152479582289
$ws.Range("O10").Value = 21.14199546738764
$ws.Range("B11").Value = 11.19068458618787
$ws.Range("C11").Value = 6.670935294790323
$ws.Range("D11").Value = 4.907719757319902
$ws.Range("E11").Value = 20.74823905097954
$ws.Range("F11").Value = 24.0041689878088
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("K11").Value = 9.964223287127377
$ws.Range("N11").Value = 17.60810499467772
$ws.Range("O11").Value = 21.13112506693341
$ws.Range("B12").Value = 11.28231972741806
$ws.Range("C12").Value = 6.718419679285143
$ws.Range("D12").Value = 4.920733082834631
$ws.Range("E12").Value = 21.00795447534692
$ws.Range("F12").Value = 24.02283910049953
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("K12").Value = 10.03208381471484
$ws.Range("N12").Value = 17.59191435628163
$ws.Range("O12").Value = 21.12780980482027
$ws.Range("B13").Value = 11.26263837180483
$ws.Range("C13").Value = 6.70822744736738
$ws.Range("D13").Value = 4.917936350733389
$ws.Range("E13").Value = 20.95228727918808
$ws.Range("F13").Value = 24.01877283906882
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("K13").Value = 10.01749928704482
$ws.Range("N13").Value = 17.59539013540391
$ws.Range("O13").Value = 21.12848814178403
$ws.Range("B14").Value = 11.19824765026313
$ws.Range("C14").Value = 6.6748573234785
$ws.Range("D14").Value = 4.908793032172591
$ws.Range("E14").Value = 20.76972643389977
$ws.Range("F14").Value = 24.00568455397814
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("K14").Value = 9.969819847312991
$ws.Range("N14").Value = 17.60676794824084
$ws.Range("O14").Value = 21.13083625032028
$ws.Range("B15").Value = 11.15864988043739
$ws.Range("C15").Value = 6.654316871884269
$ws.Range("D15").Value = 4.903175233318338
$ws.Range("E15").Value = 20.6571197878965
$ws.Range("F15").Value = 23.99780046044821
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("K15").Value = 9.940526683693861
$ws.Range("N15").Value = 17.61376990060845
$ws.Range("O15").Value = 21.13237892883348
$ws.Range("B16").Value = 10.92938846894414
$ws.Range("C16").Value = 6.535088054112456
$ws.Range("D16").Value = 4.870721172669731
$ws.Range("E16").Value = 19.99979308593342
$ws.Range("F16").Value = 23.95466303099245
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("K16").Value = 9.771355254427661
$ws.Range("N16").Value = 17.65439764702244
$ws.Range("O16").Value = 21.14281785709153
$ws.Range("B17").Value = 10.78671904719312
$ws.Range("C17").Value = 6.460608579126189
$ws.Range("D17").Value = 4.850583079403894
$ws.Range("E17").Value = 19.58578895526622
$ws.Range("F17").Value = 23.9300296423322
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("K17").Value = 9.666462744038165
$ws.Range("N17").Value = 17.67977081078456
$ws.Range("O17").Value = 21.15064624736786
$ws.Range("B18").Value = 10.70394231535827
$ws.Range("C18").Value = 6.417286972319839
$ws.Range("D18").Value = 4.838917812739639
$ws.Range("E18").Value = 19.34370152511957
$ws.Range("F18").Value = 23.91653625420192
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("K18").Value = 9.60574442494932
$ws.Range("N18").Value = 17.69453034466046
$ws.Range("O18").Value = 21.15567160565891
$ws.Range("B19").Value = 10.67579506837988
$ws.Range("C19").Value = 6.402536876182776
$ws.Range("D19").Value = 4.834954163813483
$ws.Range("E19").Value = 19.26105318901951
$ws.Range("F19").Value = 23.91208379674718
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("K19").Value = 9.585121895948403
$ws.Range("N19").Value = 17.69955614485031
$ws.Range("O19").Value = 21.15746279444247
$ws.Range("B20").Value = 10.80198129437671
$ws.Range("C20").Value = 6.468587223238588
$ws.Range("D20").Value = 4.852735375722834
$ws.Range("E20").Value = 19.6302703391353
$ws.Range("F20").Value = 23.93258209439221
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("K20").Value = 9.677669271797603
$ws.Range("N20").Value = 17.67705266877796
$ws.Range("O20").Value = 21.14975878788651
$ws.Range("B21").Value = 11.21719353416721
$ws.Range("C21").Value = 6.68467987794931
$ws.Range("D21").Value = 4.911482250621551
$ws.Range("E21").Value = 20.82351206769257
$ws.Range("F21").Value = 24.00950122961286
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("K21").Value = 9.983842930206169
$ws.Range("N21").Value = 17.60341919368833
$ws.Range("O21").Value = 21.13012479455301
$ws.Range("B22").Value = 11.48162082320528
$ws.Range("C22").Value = 6.821441630808012
$ws.Range("D22").Value = 4.949108926042308
$ws.Range("E22").Value = 21.5682972297528
$ws.Range("F22").Value = 24.06572597386447
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("K22").Value = 10.18005722045405
$ws.Range("N22").Value = 17.55676077797477
$ws.Range("O22").Value = 21.1219632088415
$ws.Range("B23").Value = 11.34115068364959
$ws.Range("C23").Value = 6.74886555781921
$ws.Range("D23").Value = 4.929098745178367
$ws.Range("E23").Value = 21.17398920530033
$ws.Range("F23").Value = 24.03517610779967
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("K23").Value = 10.07570997649765
$ws.Range("N23").Value = 17.58152960935327
$ws.Range("O23").Value = 21.12589117623879
$ws.Range("B24").Value = 10.79508357511707
$ws.Range("C24").Value = 6.46498164112905
$ws.Range("D24").Value = 4.851762595033914
$ws.Range("E24").Value = 19.61017298865642
$ws.Range("F24").Value = 23.93142604726075
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("K24").Value = 9.672604085697596
$ws.Range("N24").Value = 17.67828100435702
$ws.Range("O24").Value = 21.15015837398285
$ws.Range("B25").Value = 10.17907338109204
$ws.Range("C25").Value = 6.140429222239264
$ws.Range("D25").Value = 4.76521555593911
$ws.Range("E25").Value = 17.77168374368506
$ws.Range("F25").Value = 23.84534334191644
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("K25").Value = 9.223315184471433
$ws.Range("N25").Value = 17.78897847565323
$ws.Range("O25").Value = 21.19630335628161
